$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Turn the single run "Addig data pushing in main branch" into three runs:
#    "Addi" | "n" | "g data pushing in main branch"  (i.e. "Addig" -> "Adding",
#    spelled out as a literal run split instead of an in place text edit so
#    the run boundaries survive the save).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2Start = $p2.Range.Start

# Split right after "Addi" (4 characters into the paragraph).
$splitPoint = $d.Range($p2Start + 4, $p2Start + 4)
$splitPoint.InsertParagraphAfter()

# Paragraph 2 is now "Addi", paragraph 3 is "g data pushing in main branch"
# (still carrying the trailing bookmark). Insert a brand new paragraph
# holding just "n" between them.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.InsertAfter("n")

# Re-join "Addi" + "n" by deleting the paragraph mark between them - this
# keeps them as two distinct <w:r> runs instead of merging the text.
$p2 = $d.Paragraphs(2)
$mark = $p2.Range.End
$d.Range($mark - 1, $mark).Delete()

# Re-join the result with "g data pushing in main branch" the same way.
$p2 = $d.Paragraphs(2)
$mark = $p2.Range.End
$d.Range($mark - 1, $mark).Delete()

# ---------------------------------------------------------------------------
# 2) Move the _GoBack bookmark off this paragraph: it will reappear at the
#    end of the new "Adding new line" paragraph below.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 3) Add the "Adding new line" paragraph right after the edited paragraph,
#    with a throwaway tail of filler text. Placing the bookmark while real
#    text still follows it keeps the engine from mis-resolving a bookmark
#    that would otherwise sit exactly at the document's trailing edge.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$newParaStart = $p3.Range.Start
$p3.Range.InsertAfter("Adding new lineFILLER")

$boundary = $newParaStart + ("Adding new line".Length)
$d.Bookmarks.Add("_GoBack", $d.Range($boundary, $boundary))

# Remove the filler text now that the bookmark is anchored in place.
$d.Range($boundary, $boundary + 6).Delete()

# ---------------------------------------------------------------------------
# 4) Append the trailing empty paragraph, once everything else is settled.
# ---------------------------------------------------------------------------
$d.Content.InsertParagraphAfter()
